# Upload new version with timestamp
# ---------------------------------------------------------------------------
# The underlying report re-pulled stock data: 3 new drugs were added to the
# (alphabetically sorted) product list - COLDATREXY, DELTARHINO and HIBIOTIC -
# and every row's "current balance" / "sale price" / "transactions" figures
# were refreshed. The grand total and the two footer rows (page no. /
# developer credit) simply slide three rows further down the sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the three new product rows right before the totals row
#    (old row 18). Excel shifts everything below - totals, footer rows and
#    their merged ranges - down by three automatically.
$ws.Rows("18:20").Insert()

# 2) The newly inserted rows come back blank/unformatted, so clone the exact
#    look (fonts, fills, borders, number formats) of the last product row
#    (row 17) onto each of them.
$ws.Range("A17:N17").Copy()
$ws.Range("A18:N18").PasteSpecial(-4122)
$ws.Range("A17:N17").Copy()
$ws.Range("A19:N19").PasteSpecial(-4122)
$ws.Range("A17:N17").Copy()
$ws.Range("A20:N20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Restore the alternating row heights used throughout the product table.
$ws.Rows("18").RowHeight = 25.5
$ws.Rows("19").RowHeight = 24.75
$ws.Rows("20").RowHeight = 25.5

# PasteSpecial(formats) does not carry merged-cell state, so re-merge the
# same three column groups used by every other product row.
$ws.Range("B18:G18").Merge()
$ws.Range("H18:K18").Merge()
$ws.Range("L18:M18").Merge()
$ws.Range("B19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("B20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()

function Set-RowData {
    param($r, $a, $b, $h, $l, $n)
    $ws.Cells.Item($r, 1).Value = $a       # A - م (sequence number)
    $ws.Cells.Item($r, 2).Value = $b       # B - الاسم (product name)
    $ws.Cells.Item($r, 8).Value = $h       # H - الرصيد الحالي
    $ws.Cells.Item($r, 12).Value = $l      # L - سعر البيع
    $ws.Cells.Item($r, 14).Value = $n      # N - عدد التعاملات
}

# 4) Refill the whole product table (rows 4-20) with the refreshed data -
#    now 17 rows instead of 14, still sorted A->Z by product name.
Set-RowData 4 1 "AVEROTHIAZIDE 5/20/12.5MG 30 F.C. TABS" "1:0" 93 "1:0"
Set-RowData 5 2 "CIPRODIAZOLE 500/500MG 20 F.C. TABS" "1:0" 74 "0:2"
Set-RowData 6 3 "CLAVIMOX 1 GM 12 F.C.TABS." "2:1" 43.33 "0:0"
Set-RowData 7 4 "COLDATREXY 30 F.C. TABS" "1:1" 21 "0:0"
Set-RowData 8 5 "CONVENTIN XR 600MG 30 TABS." "0:2" 0 "0:0"
Set-RowData 9 6 "DELTARHINO NASAL SPRAY 15 ML" "0:0" 30 "1:0"
Set-RowData 10 7 "GAST-REG 200 MG 30 TABS." "1:3" 28 "0:0"
Set-RowData 11 8 "HELI-CURE 14 ENTERIC COATED TAB" "2:1" 120 "0:2"
Set-RowData 12 9 "HIBIOTIC 625MG 16 TAB." "1:1" 71.5 "0:2"
Set-RowData 13 10 "OTRIVIN 0.1% ADULT NASAL DROPS 15 ML" "6:0" 24 "1:0"
Set-RowData 14 11 "SENSODERM صابون مرطب" "0:0" 50 "1:0"
Set-RowData 15 12 "VASTAFLAM 50MG 20 SUGAR COATED TAB." "2:1" 18 "0:2"
Set-RowData 16 13 "VOLTAREN 75MG/3ML 3 AMP." "5:3" 34 "0:3"
Set-RowData 17 14 "XILOPRED 16MG  20TAB" "0:1" 37 "0:2"
Set-RowData 18 15 "سرنجات 3 سم" "-2:0" 4 "2:0"
Set-RowData 19 16 "كالونا " "-1:0" 15 "1:0"
Set-RowData 20 17 "مناديل FINE" "14:0" 30 "1:0"

# 5) Grand total (K21, merged K21:N21) now sums the 17 refreshed prices.
$ws.Range("K21").Value = 692.83

# 6) Minor autofit nudge on the developer-credit footer row (now row 22)
#    that Excel applies whenever the sheet is regenerated.
$ws.Rows("22").RowHeight = 16.5

Write-Host "Edit applied"
